$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 0.1818801724491279
$ws.Range("J2").Value = 0.2500781318045117
$ws.Range("M2").Value = 0.4260053333333333
$ws.Range("N2").Value = 1.278016
$ws.Range("O2").Value = 0.02405532912416773
$ws.Range("P2").Value = 0.02531756756689831
$ws.Range("Q2").Value = 0.01386249755022222
$ws.Range("R2").Value = 0.124762477952
$ws.Range("S2").Value = 0.004375187409424155
$ws.Range("T2").Value = 0.006331369998964425

# Row 3
$ws.Range("I3").Value = 0.1818801724491279
$ws.Range("J3").Value = 0.2500781318045117
$ws.Range("O3").Value = 0.05879323641880037
$ws.Range("P3").Value = 0.06187825274916518
$ws.Range("S3").Value = 0.01069332397869376
$ws.Range("T3").Value = 0.01547439784683862

# Row 4
$ws.Range("I4").Value = 0.1818801724491279
$ws.Range("J4").Value = 0.2500781318045117
$ws.Range("M4").Value = 5.850740666666667
$ws.Range("N4").Value = 17.552222
$ws.Range("O4").Value = 0.3303749538898241
$ws.Range("P4").Value = 0.3477104875323931
$ws.Range("Q4").Value = 0.1903870017871111
$ws.Range("R4").Value = 1.713483016084
$ws.Range("S4").Value = 0.06008865358635389
$ws.Range("T4").Value = 0.08695478913093682

# Row 5
$ws.Range("I5").Value = 0.1818801724491279
$ws.Range("J5").Value = 0.2500781318045117
$ws.Range("M5").Value = 2.648771
$ws.Range("N5").Value = 5.297542
$ws.Range("O5").Value = 0.1495686865725097
$ws.Range("P5").Value = 0.1049445996947469
$ws.Range("Q5").Value = 0.08619277418733334
$ws.Range("R5").Value = 0.5171566451239999
$ws.Range("S5").Value = 0.02720357850679763
$ws.Range("T5").Value = 0.02624434943463462

# Row 6
$ws.Range("I6").Value = 0.1818801724491279
$ws.Range("J6").Value = 0.2500781318045117
$ws.Range("M6").Value = 7.742685666666667
$ws.Range("N6").Value = 23.228057
$ws.Range("O6").Value = 0.4372077939946981
$ws.Range("P6").Value = 0.4601490924567965
$ws.Range("Q6").Value = 0.2519521533837778
$ws.Range("R6").Value = 2.267569380454
$ws.Range("S6").Value = 0.07951942896785848
$ws.Range("T6").Value = 0.1150732253931372

# Row 7
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.5
$ws.Range("G7").Value = 0.146372
$ws.Range("H7").Value = 0.292744
$ws.Range("I7").Value = 0.8181198275508721
$ws.Range("J7").Value = 0.7499218681954883
$ws.Range("M7").Value = 0.4260053333333333
$ws.Range("N7").Value = 1.278016
$ws.Range("O7").Value = 0.02405532912416773
$ws.Range("P7").Value = 0.02531756756689831
$ws.Range("Q7").Value = 0.06235525265066667
$ws.Range("R7").Value = 0.374131515904
$ws.Range("S7").Value = 0.01968014171474357
$ws.Range("T7").Value = 0.01898619756793388

# Row 8
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.5
$ws.Range("G8").Value = 0.146372
$ws.Range("H8").Value = 0.292744
$ws.Range("I8").Value = 0.8181198275508721
$ws.Range("J8").Value = 0.7499218681954883
$ws.Range("O8").Value = 0.05879323641880037
$ws.Range("P8").Value = 0.06187825274916518
$ws.Range("Q8").Value = 0.1524014530053333
$ws.Range("R8").Value = 0.9144087180320001
$ws.Range("S8").Value = 0.04809991244010661
$ws.Range("T8").Value = 0.04640385490232657

# Row 9
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.5
$ws.Range("G9").Value = 0.146372
$ws.Range("H9").Value = 0.292744
$ws.Range("I9").Value = 0.8181198275508721
$ws.Range("J9").Value = 0.7499218681954883
$ws.Range("M9").Value = 5.850740666666667
$ws.Range("N9").Value = 17.552222
$ws.Range("O9").Value = 0.3303749538898241
$ws.Range("P9").Value = 0.3477104875323931
$ws.Range("Q9").Value = 0.8563846128613334
$ws.Range("R9").Value = 5.138307677168
$ws.Range("S9").Value = 0.2702863003034702
$ws.Range("T9").Value = 0.2607556984014563

# Row 10
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.5
$ws.Range("G10").Value = 0.146372
$ws.Range("H10").Value = 0.292744
$ws.Range("I10").Value = 0.8181198275508721
$ws.Range("J10").Value = 0.7499218681954883
$ws.Range("M10").Value = 2.648771
$ws.Range("N10").Value = 5.297542
$ws.Range("O10").Value = 0.1495686865725097
$ws.Range("P10").Value = 0.1049445996947469
$ws.Range("Q10").Value = 0.387705908812
$ws.Range("R10").Value = 1.550823635248
$ws.Range("S10").Value = 0.1223651080657121
$ws.Range("T10").Value = 0.07870025026011224

# Row 11
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0.5
$ws.Range("G11").Value = 0.146372
$ws.Range("H11").Value = 0.292744
$ws.Range("I11").Value = 0.8181198275508721
$ws.Range("J11").Value = 0.7499218681954883
$ws.Range("M11").Value = 7.742685666666667
$ws.Range("N11").Value = 23.228057
$ws.Range("O11").Value = 0.4372077939946981
$ws.Range("P11").Value = 0.4601490924567965
$ws.Range("Q11").Value = 1.133312386401333
$ws.Range("R11").Value = 6.799874318408
$ws.Range("S11").Value = 0.3576883650268396
$ws.Range("T11").Value = 0.3450758670636593
